# Insert a new data row at row 20 (pushing the existing rows 20-74 down to
# 21-75, which also grows the used range from A1:R74 to A1:R75), then
# populate the newly inserted row with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(20).Insert()

$ws.Cells.Item(20,1).Value  = 1
$ws.Cells.Item(20,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(20,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(20,4).Value  = 44624
$ws.Cells.Item(20,5).Value  = 15
$ws.Cells.Item(20,6).Value  = 100112038
$ws.Cells.Item(20,7).Value  = "Cebollín baby"
$ws.Cells.Item(20,8).Value  = "Sin especificar"
$ws.Cells.Item(20,9).Value  = "Primera"
$ws.Cells.Item(20,10).Value = 250
$ws.Cells.Item(20,11).Value = 2800
$ws.Cells.Item(20,12).Value = 3000
$ws.Cells.Item(20,13).Value = 2900
$ws.Cells.Item(20,14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(20,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(20,16).Value = 1450
$ws.Cells.Item(20,17).Value = 2
$ws.Cells.Item(20,18).Value = "Hortaliza"
